$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = "{value_covn1_0|<ND>|<MISSING>}"
$ws.Range("G3").Value = "{value_covn1_1|<ND>|<MISSING>}"
$ws.Range("H3").Value = "{value_covn1_2|<ND>|<MISSING>}"

$ws.Range("V3").Value = "{value_d3l_del_0|<ND>|<MISSING>}"
$ws.Range("W3").Value = "{value_d3l_del_1|<ND>|<MISSING>}"
$ws.Range("X3").Value = "{value_d3l_del_2|<ND>|<MISSING>}"

$ws.Range("Y3").Select()
